# Automatische test-sync: 2025-08-14 22:11:50
# Appends the new mail-log entry (row 40) to the "Logs" sheet, extends the
# conditional-formatting ranges that were anchored at row 39 so they cover
# the new row, and bumps the "Intern verzoek / Actie voor medewerker"
# tally on the "Dashboard" sheet from 31 to 32.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Logs sheet: append row 40
# ---------------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A40").Value = "Vraag over product"
$logs.Range("B40").Value = "documentatie@testbedrijf123.nl"
$logs.Range("C40").Value = "Is de EcoPro-700 nog op voorraad?"
$logs.Range("D40").Value = "Intern verzoek / Actie voor medewerker"
$logs.Range("E40").Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@testbedrijf123.nl."
$logs.Range("F40").Value = "2025-08-14 22:10:52"
$logs.Range("G40").Value = "Nee"
$logs.Range("H40").Value = "Ja"
$logs.Range("I40").Value = "Nee"
$logs.Range("J40").Value = "Nee"

# ---------------------------------------------------------------------
# 2. Logs sheet: stretch the conditional formatting sqref from row 39 to
#    row 40 for every column that had it (D, G, H, I, J). Each of those
#    conditionalFormatting blocks shares a single sqref across all of its
#    cfRules, so re-pointing rule 1's AppliesTo range re-points the whole
#    block.
# ---------------------------------------------------------------------
$logs.Range("D2:D39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("D2:D40"))
$logs.Range("G2:G39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("G2:G40"))
$logs.Range("H2:H39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("H2:H40"))
$logs.Range("I2:I39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("I2:I40"))
$logs.Range("J2:J39").FormatConditions.Item(1).ModifyAppliesToRange($logs.Range("J2:J40"))

# ---------------------------------------------------------------------
# 3. Dashboard sheet: bump the count for "Intern verzoek / Actie voor
#    medewerker" from 31 to 32.
# ---------------------------------------------------------------------
$dashboard = $wb.Worksheets.Item("Dashboard")
$dashboard.Range("B2").Value = 32
